$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDiCECpDoC")
$ws.Range("B1").Value = "Perc Decline per Doubling (dimensionless)"
